$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 (Camarilla pivot) updates
$ws1.Range("C3").Value = 4639.71
$ws1.Range("D3").Value = 197.36
$ws1.Range("E3").Value = 46796.15
$ws1.Range("F3").Value = 67463.37
$ws1.Range("G3").Value = 693.75
$ws1.Range("H3").Value = 226.76
$ws1.Range("I3").Value = 165.46
$ws1.Range("J3").Value = 1247.72
$ws1.Range("K3").Value = 185.74
$ws1.Range("C4").Value = 4602.58
$ws1.Range("D4").Value = 196.14
$ws1.Range("E4").Value = 46711.25
$ws1.Range("F4").Value = 67161.99000000001
$ws1.Range("G4").Value = 692.17
$ws1.Range("H4").Value = 225.5
$ws1.Range("I4").Value = 165.12
$ws1.Range("J4").Value = 1240.62
$ws1.Range("K4").Value = 185.12
$ws1.Range("C5").Value = 4565.45
$ws1.Range("D5").Value = 194.92
$ws1.Range("E5").Value = 46626.35
$ws1.Range("F5").Value = 66860.60000000001
$ws1.Range("G5").Value = 690.59
$ws1.Range("H5").Value = 224.23
$ws1.Range("I5").Value = 164.77
$ws1.Range("J5").Value = 1233.51
$ws1.Range("K5").Value = 184.51
$ws1.Range("C6").Value = 4521.73
$ws1.Range("D6").Value = 193.46
$ws1.Range("E6").Value = 46522.68
$ws1.Range("F6").Value = 66494.3
$ws1.Range("G6").Value = 688.7
$ws1.Range("H6").Value = 222.69
$ws1.Range("I6").Value = 164.36
$ws1.Range("J6").Value = 1224.91
$ws1.Range("K6").Value = 183.78
$ws1.Range("C7").Value = 4507.1
$ws1.Range("D7").Value = 192.97
$ws1.Range("E7").Value = 46487.99
$ws1.Range("F7").Value = 66371.75999999999
$ws1.Range("G7").Value = 688.0599999999999
$ws1.Range("H7").Value = 222.17
$ws1.Range("I7").Value = 164.22
$ws1.Range("J7").Value = 1222.03
$ws1.Range("K7").Value = 183.53
$ws1.Range("C8").Value = 4492.56
$ws1.Range("D8").Value = 192.49
$ws1.Range("E8").Value = 46453.53
$ws1.Range("F8").Value = 66250.00999999999
$ws1.Range("G8").Value = 687.4299999999999
$ws1.Range("H8").Value = 221.66
$ws1.Range("I8").Value = 164.09
$ws1.Range("J8").Value = 1219.17
$ws1.Range("K8").Value = 183.29
$ws1.Range("C9").Value = 4463.44
$ws1.Range("D9").Value = 191.51
$ws1.Range("E9").Value = 46384.47
$ws1.Range("F9").Value = 66005.99000000001
$ws1.Range("G9").Value = 686.17
$ws1.Range("H9").Value = 220.64
$ws1.Range("I9").Value = 163.81
$ws1.Range("J9").Value = 1213.43
$ws1.Range("K9").Value = 182.81
$ws1.Range("C10").Value = 4448.9
$ws1.Range("D10").Value = 191.03
$ws1.Range("E10").Value = 46350.01
$ws1.Range("F10").Value = 65884.24000000001
$ws1.Range("G10").Value = 685.54
$ws1.Range("H10").Value = 220.13
$ws1.Range("I10").Value = 163.68
$ws1.Range("J10").Value = 1210.57
$ws1.Range("K10").Value = 182.57
$ws1.Range("C11").Value = 4434.27
$ws1.Range("D11").Value = 190.54
$ws1.Range("E11").Value = 46315.32
$ws1.Range("F11").Value = 65761.7
$ws1.Range("G11").Value = 684.9
$ws1.Range("H11").Value = 219.61
$ws1.Range("I11").Value = 163.54
$ws1.Range("J11").Value = 1207.69
$ws1.Range("K11").Value = 182.32
$ws1.Range("C12").Value = 4390.55
$ws1.Range("D12").Value = 189.08
$ws1.Range("E12").Value = 46211.65
$ws1.Range("F12").Value = 65395.4
$ws1.Range("G12").Value = 683
$ws1.Range("H12").Value = 218.07
$ws1.Range("I12").Value = 163.12
$ws1.Range("J12").Value = 1199.08
$ws1.Range("K12").Value = 181.59
$ws1.Range("C13").Value = 4353.42
$ws1.Range("D13").Value = 187.86
$ws1.Range("E13").Value = 46126.75
$ws1.Range("F13").Value = 65094.01
$ws1.Range("G13").Value = 681.4299999999999
$ws1.Range("H13").Value = 216.8
$ws1.Range("I13").Value = 162.78
$ws1.Range("J13").Value = 1191.98
$ws1.Range("K13").Value = 180.98
$ws1.Range("C14").Value = 4316.29
$ws1.Range("D14").Value = 186.64
$ws1.Range("E14").Value = 46041.85
$ws1.Range("F14").Value = 64792.63
$ws1.Range("G14").Value = 679.85
$ws1.Range("H14").Value = 215.54
$ws1.Range("I14").Value = 162.44
$ws1.Range("J14").Value = 1184.88
$ws1.Range("K14").Value = 180.36

# Sheet2 (raw OHLC) updates
$ws2.Range("B2").Value = "19APR2021"
$ws2.Range("C2").Value = 4442
$ws2.Range("D2").Value = 4562
$ws2.Range("E2").Value = 4403
$ws2.Range("F2").Value = 4478
$ws2.Range("G2").Value = 4434
$ws2.Range("B3").Value = "27APR2021"
$ws2.Range("C3").Value = 190
$ws2.Range("D3").Value = 195.3
$ws2.Range("E3").Value = 190
$ws2.Range("F3").Value = 192
$ws2.Range("G3").Value = 189.1
$ws2.Range("B4").Value = "04JUN2021"
$ws2.Range("C4").Value = 46545
$ws2.Range("D4").Value = 46777
$ws2.Range("E4").Value = 46400
$ws2.Range("F4").Value = 46419
$ws2.Range("G4").Value = 46593
$ws2.Range("C5").Value = 66786
$ws2.Range("D5").Value = 67293
$ws2.Range("E5").Value = 65961
$ws2.Range("F5").Value = 66128
$ws2.Range("G5").Value = 66983
$ws2.Range("B6").Value = "30APR2021"
$ws2.Range("C6").Value = 686
$ws2.Range("D6").Value = 688.65
$ws2.Range("E6").Value = 681.75
$ws2.Range("F6").Value = 686.8
$ws2.Range("G6").Value = 689.85
$ws2.Range("B7").Value = "30APR2021"
$ws2.Range("C7").Value = 226.3
$ws2.Range("D7").Value = 226.3
$ws2.Range("E7").Value = 220.7
$ws2.Range("F7").Value = 221.15
$ws2.Range("G7").Value = 227.7
$ws2.Range("B8").Value = "30APR2021"
$ws2.Range("C8").Value = 164
$ws2.Range("D8").Value = 164.4
$ws2.Range("E8").Value = 162.9
$ws2.Range("F8").Value = 163.95
$ws2.Range("G8").Value = 164.45
$ws2.Range("B9").Value = "30APR2021"
$ws2.Range("C9").Value = 1243
$ws2.Range("D9").Value = 1243
$ws2.Range("E9").Value = 1211.7
$ws2.Range("F9").Value = 1216.3
$ws2.Range("G9").Value = 1252.7
$ws2.Range("B10").Value = "30APR2021"
$ws2.Range("C10").Value = 181.5
$ws2.Range("D10").Value = 183.3
$ws2.Range("E10").Value = 180.65
$ws2.Range("F10").Value = 183.05
$ws2.Range("G10").Value = 182.25
